# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the upstream data source.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 387
$wsExhibit.Range("F4").Value = 4940
$wsExhibit.Range("F5").Value = 28
$wsExhibit.Range("F6").Value = 28
$wsExhibit.Range("F8").Value = 492

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 387
$wsAll.Range("F4").Value = 4940
$wsAll.Range("F6").Value = 28
$wsAll.Range("F7").Value = 28
$wsAll.Range("F10").Value = 492
